$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1. Insert a new column at Y ("text" function-category) for the new
#    "outputToCloud(resource)" / spellcheck-related function list.
#    This pushes the existing web/webalert/webcookie/ws/ws.async/xml columns
#    from Y..AD to Z..AE.
# ---------------------------------------------------------------------------
$ws.Columns("Y").Insert()
$ws.Cells.Item(1, 25).Value2 = "text"
$ws.Cells.Item(2, 25).Value2 = "spellCheck(var,profile,text)"

# ---------------------------------------------------------------------------
# 2. Add "text" to the "target" list (column A), keeping it alphabetically
#    sorted between "step" and "web" (row 25), shifting the rest down by one
#    row. A plain Range.Insert() on this runtime shifts the whole row across
#    all columns, so shift only column A manually (bottom-up to avoid
#    clobbering).
# ---------------------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $v = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 1).Value2 = $v
}
$ws.Cells.Item(25, 1).Value2 = "text"

# ---------------------------------------------------------------------------
# 3. Add "outputToCloud(resource)" to the "base" function list (column E),
#    alphabetically between "macro(file,sheet,name)" and
#    "prependText(var,prependWith)" (row 22), shifting the rest down by one
#    row (column E only).
# ---------------------------------------------------------------------------
for ($r = 38; $r -ge 22; $r--) {
    $v = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r + 1, 5).Value2 = $v
}
$ws.Cells.Item(22, 5).Value2 = "outputToCloud(resource)"

# ---------------------------------------------------------------------------
# 4. Fix up the defined names that refer to ranges affected by the shifts
#    above, and register the new "text" named range.
# ---------------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
